$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Insert a brand-new step "0: Go to ItemColumns and add/remove the new
# column" as its own paragraph right before the existing
# "1: Go to InventarioManager ..." paragraph.
$target1 = $d.Paragraphs.Item(2)
$target1.Range.InsertParagraphBefore()
$d.Paragraphs.Item(2).Range.Text = "0: Go to ItemColumns and add/remove the new column"

# --- Change 2 -----------------------------------------------------------
# Right after the "2: Go to the part where it gets the JSON value and
# add/remove a new row;" paragraph, add two new step paragraphs (3 and 4)
# followed by a brand-new blank paragraph (the pre-existing blank
# paragraph that used to follow step 2 is left untouched, just pushed
# further down).
$d.Content.Find.Execute(
    "2: Go to the part where it gets the JSON value and add/remove a new row;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2: Go to the part where it gets the JSON value and add/remove a new row;^p" + `
    "3: Go to CreateAddItemForm and add/remove a field to the corresponding function^p" + `
    "4: Go to the corresponding AddItem php file and add/remove the field^p",
    2)
